# Read basic auth from run manager for book tests
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

# Row 4 corresponds to the "BookTests" run config.
# Previously the auth info was stuffed into the "headers" column (E4) as a
# pre-built Basic auth header. Switch it to use the authType/authParams
# columns (B4/C4) instead, and clear the old headers value.
$ws.Range("B4").Value = "basic"
$ws.Range("C4").Value = "bWVyX2dyZXk=:U2hlcGhlcmRAMTIz"
$ws.Range("E4").Value = ""
